$wb = $excel.ActiveWorkbook

$wsFolderPerms = $wb.Worksheets.Item("Folder Perms")
$wsPermsRef = $wb.Worksheets.Item("Perms Reference")
$wsServerList = $wb.Worksheets.Item("Server List")

# Row 4 on "Folder Perms": the "new-folder" row for SPORe (E:\WebTools\PO\),
# written first so the new shared strings are appended in the same order as
# the target workbook (new-folder, E:\WebTools\PO\, ..., E:\WebTools\, n/a).
$wsFolderPerms.Range("A4").Value = "Web"
$wsFolderPerms.Range("D4").Value = "new-folder"
$wsFolderPerms.Range("B4").Value = "E:\WebTools\PO\"
$wsFolderPerms.Range("E4").Value = "y"

# Document the new "new-folder" perm value on the reference sheet
$wsPermsRef.Range("A27").Value = "new-folder"
$wsPermsRef.Range("B27").Value = "this tells processing script to just create the folder"

# Row 3 on "Folder Perms": the regular Modify row for the parent E:\WebTools\ folder
$wsFolderPerms.Range("A3").Value = "Web"
$wsFolderPerms.Range("B3").Value = "E:\WebTools\"
$wsFolderPerms.Range("C3").Value = "IIS_IUSRS"
$wsFolderPerms.Range("D3").Value = "Modify"
$wsFolderPerms.Range("E3").Value = "y"

# Last new value
$wsFolderPerms.Range("C4").Value = "n/a"

# Update sheet selections to reflect where the edits were made
$wsFolderPerms.Activate()
$wsFolderPerms.Range("A2:E4").Select()

$wsPermsRef.Activate()
$wsPermsRef.Range("A28").Select()

# Restore the workbook's original active sheet/selection
$wsServerList.Activate()
$wsServerList.Range("A10").Select()
